$wb = $excel.ActiveWorkbook

# Sheet 1: Accounts
$wsAccounts = $wb.Worksheets.Item("Accounts")
$wsAccounts.Range("C2").Value = 1990
$wsAccounts.Range("C4").Value = 205

# Sheet 2: Sales
$wsSales = $wb.Worksheets.Item("Sales")
$wsSales.Range("B2").Value = 510
$wsSales.Range("B3").Value = 2985
$wsSales.Range("B4").Value = 430
